# Update daily power records
# Row 56 (2018-10-02, serial 43380): fill in Start Time / End Time that were
# previously left at 0, which also recalculates the dependent Duration /
# Second Duration / Absolute Value formulas for that row.
# Row 57 (2018-10-03, serial 43381): fill in the Start Time (B57), which was
# previously blank.
# Finally move the active selection to C57 to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 56 - Start Time 20:30:00, End Time 23:59:00
$ws.Range("B56").Value = 0.85416666666666663
$ws.Range("C56").Value = 0.99930555555555556

# Row 57 - Start Time 00:00:00 (newly entered)
$ws.Range("B57").Value = 0

# Move selection to C57
$ws.Range("C57").Select()
